$wb = $excel.ActiveWorkbook

# --- Sheet "Orders" (sheet1): replace order rows 2-27 with the new 7-row order list ---
$ws = $wb.Worksheets.Item("Orders")

# Clear out all the old data below the header (rows 2-27, columns A-L) without
# shifting cells, so the used range collapses back down once the new, shorter
# data block is written.
$ws.Range("A2:L27").ClearContents()

# The PackageID (col A) and Number (col F) columns hold numeric-looking values
# that must stay stored as TEXT (as in the source workbook). Force a text
# number format on exactly the cells that get a value, so no stray empty
# cells are introduced on the rows that have no PackageID.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F8").NumberFormat = "@"

# New order data (PackageID in col A, FlowerName in col C, Number in col F)
$ws.Range("A2").Value = "35"
$ws.Range("C2").Value = "111_绣球单瓣紫粉_Hydrangea Purple&Pink S_Hydrangea L._1stem"
$ws.Range("F2").Value = "20"

$ws.Range("C3").Value = "771_美洲茶_undefined_undefined_1bunch"
$ws.Range("F3").Value = "15"

$ws.Range("C4").Value = "495_大飞燕深粉色_delphinium pink_undefined_1bunch"
$ws.Range("F4").Value = "21"

$ws.Range("C5").Value = "457_茴香花_lace flower yellow_undefined_1bunch"
$ws.Range("F5").Value = "20"

$ws.Range("C6").Value = "455_粉星花_tweedia pink_undefined_1bunch"
$ws.Range("F6").Value = "10"

$ws.Range("C7").Value = "456_蕾丝白色_lace flower white_undefined_1bunch"
$ws.Range("F7").Value = "20"

$ws.Range("C8").Value = "574_迷你菊白_undefined_undefined_1bunch"
$ws.Range("F8").Value = "35"

# --- Sheet "Summary" (sheet2): TotalNumber string reflects the new Number column ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("G2").NumberFormat = "@"
$ws2.Range("G2").Value = "020152120102035"
